# order_matrix_12_A_block3_VR.xlsx edit:
#   - participant id changes from 12 to 14 throughout (col A)
#   - every existing block_num (col I) drops from 4 to 2
#   - instruction/video paths (col H) move under an /VR/ subfolder (and a
#     few are swapped for different clips), video_id/dimension (cols E/F)
#     change for the "video" rows, and the last four rows (14-17) pick up
#     new descriptions/paths that used to belong to a 4th (luminance) block
#   - four brand new rows (18-21) are appended, carrying over what used to
#     be the luminance sub-block (now suprablock_count 3, block_order 2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (participant) holds numeric-looking text ("14") but must stay
# a text cell, matching the workbook's original inlineStr/string typing -
# so force a text number format before writing into every row touched
# below. Column E (video_id) needs the same treatment, but only on the
# specific rows that get a video_id value written.
$ws.Range("A2:A21").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "14"
$ws.Range("H2").Value = "./instructions_videos/VR/block_2_text.mp4"
$ws.Range("I2").Value = 2

# Row 3
$ws.Range("A3").Value = "14"
$ws.Range("E3").Value = "1"
$ws.Range("F3").Value = "arousal"
$ws.Range("G3").Value = "inverse "
$ws.Range("H3").Value = "../stimuli/exp_videos/VR/1.mp4"
$ws.Range("I3").Value = 2

# Row 4
$ws.Range("A4").Value = "14"
$ws.Range("H4").Value = "./instructions_videos/VR/post_stimulus_self_report.mp4"
$ws.Range("I4").Value = 2

# Row 5
$ws.Range("A5").Value = "14"
$ws.Range("H5").Value = "./instructions_videos/VR/mareo.mp4"
$ws.Range("I5").Value = 2

# Row 6
$ws.Range("A6").Value = "14"
$ws.Range("H6").Value = "./instructions_videos/VR/block_2_text_reminder.mp4"
$ws.Range("I6").Value = 2

# Row 7
$ws.Range("A7").Value = "14"
$ws.Range("E7").Value = "5"
$ws.Range("F7").Value = "arousal"
$ws.Range("H7").Value = "../stimuli/exp_videos/VR/5.mp4"
$ws.Range("I7").Value = 2

# Row 8
$ws.Range("A8").Value = "14"
$ws.Range("H8").Value = "./instructions_videos/VR/post_stimulus_self_report.mp4"
$ws.Range("I8").Value = 2

# Row 9
$ws.Range("A9").Value = "14"
$ws.Range("H9").Value = "./instructions_videos/VR/mareo.mp4"
$ws.Range("I9").Value = 2

# Row 10
$ws.Range("A10").Value = "14"
$ws.Range("H10").Value = "./instructions_videos/VR/block_2_text_reminder.mp4"
$ws.Range("I10").Value = 2

# Row 11
$ws.Range("A11").Value = "14"
$ws.Range("E11").Value = "11"
$ws.Range("F11").Value = "arousal"
$ws.Range("H11").Value = "../stimuli/exp_videos/VR/11.mp4"
$ws.Range("I11").Value = 2

# Row 12
$ws.Range("A12").Value = "14"
$ws.Range("H12").Value = "./instructions_videos/VR/post_stimulus_self_report.mp4"
$ws.Range("I12").Value = 2

# Row 13
$ws.Range("A13").Value = "14"
$ws.Range("H13").Value = "./instructions_videos/VR/mareo.mp4"
$ws.Range("I13").Value = 2

# Row 14
$ws.Range("A14").Value = "14"
$ws.Range("H14").Value = "./instructions_videos/VR/block_2_text_reminder.mp4"
$ws.Range("I14").Value = 2
$ws.Range("L14").Value = "audio_instruction"

# Row 15
$ws.Range("A15").Value = "14"
$ws.Range("E15").Value = "10"
$ws.Range("F15").Value = "arousal"
$ws.Range("H15").Value = "../stimuli/exp_videos/VR/10.mp4"
$ws.Range("I15").Value = 2
$ws.Range("L15").Value = "video"

# Row 16
$ws.Range("A16").Value = "14"
$ws.Range("H16").Value = "./instructions_videos/VR/post_stimulus_self_report.mp4"
$ws.Range("I16").Value = 2
$ws.Range("L16").Value = "post_stimulus_self_report"

# Row 17
$ws.Range("A17").Value = "14"
$ws.Range("H17").Value = "./instructions_videos/VR/mareo.mp4"
$ws.Range("I17").Value = 2
$ws.Range("L17").Value = "motion_sickness"

# Row 18 (new)
$ws.Range("A18").Value = "14"
$ws.Range("B18").Value = "A_block3"
$ws.Range("C18").Value = "VR"
$ws.Range("H18").Value = "./instructions_videos/VR/luminance_instructions_direct.mp4"
$ws.Range("J18").Value = 2
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = "luminance_instructions"

# Row 19 (new)
$ws.Range("A19").Value = "14"
$ws.Range("B19").Value = "A_block3"
$ws.Range("C19").Value = "VR"
$ws.Range("D19").Value = 5
$ws.Range("F19").Value = "luminance"
$ws.Range("G19").Value = "inverse "
$ws.Range("H19").Value = "../stimuli/exp_videos/VR/green_intensity_video_7.mp4"
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = "luminance"

# Row 20 (new)
$ws.Range("A20").Value = "14"
$ws.Range("B20").Value = "A_block3"
$ws.Range("C20").Value = "VR"
$ws.Range("H20").Value = "./instructions_videos/VR/confidence_luminance_practice_instructions_text.mp4"
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = "confidence_luminance_instructions"

# Row 21 (new)
$ws.Range("A21").Value = "14"
$ws.Range("B21").Value = "A_block3"
$ws.Range("C21").Value = "VR"
$ws.Range("H21").Value = "./instructions_videos/VR/rest_suprablock_text.mp4"
$ws.Range("J21").Value = 2
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = "rest_suprablock"

Write-Host "edit applied"
